# This edit re-sorts/permutes the weekly price records (rows 2-13) of the
# "Haba" (fava bean) sheet so that the rows reflect the updated weekly
# reporting order. Only the per-record fields (Fecha, Volumen, Precio
# minimo/maximo/promedio, Origen, Precio $/Kg) move between rows; the
# market/region/category descriptive columns stay identical in every row
# so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D, J, K, L, M, O, P (row => values)
$data = @{
    2  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 }
    3  = @{ D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
    4  = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    5  = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    6  = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 }
    7  = @{ D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    8  = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    9  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
    10 = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    11 = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    12 = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    13 = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D    # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio $/Kg
}
